$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("K2").Value = 13
$ws.Range("N2").Value = 12
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 23
$ws.Range("U2").Value = 13
$ws.Range("W2").Value = 15
$ws.Range("Z2").Value = 6
$ws.Range("AB2").Value = 20

# Row 3 updates
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 6
$ws.Range("E3").Value = 12
$ws.Range("H3").Value = 8
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 17

# Row 4 updates
$ws.Range("L4").Value = 3
$ws.Range("N4").Value = 4
$ws.Range("Q4").Value = 3
$ws.Range("S4").Value = 7
$ws.Range("T4").Value = 17
$ws.Range("W4").Value = 10
$ws.Range("AA4").Value = 3
$ws.Range("AB4").Value = 24
